$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns hold numeric-looking text (e.g. "1.000", "241.19").
# Force the range to text format first so Excel doesn't auto-convert these
# assignments into real numbers (which would drop formatting like trailing
# zeros or the "x.xxx.xx" grouped-price strings).
$ws.Range("D2:E51").NumberFormat = "@"

# Row 8 / Row 9 swapped which coin they describe (OKB now ranks above Dogecoin)
$ws.Range("B8").Value = "OKB"
$ws.Range("C8").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D8").Value = "44.87"
$ws.Range("E8").Value = "  +7.30%  "

$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "0.07524"
$ws.Range("E9").Value = "  +1.18%  "

$ws.Range("D2").Value = "29.020.11"
$ws.Range("E2").Value = "  -0.55%  "

$ws.Range("D3").Value = "1.829.23"
$ws.Range("E3").Value = "  -0.33%  "

$ws.Range("D4").Value = "0.9990"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "241.19"
$ws.Range("E5").Value = "  -0.32%  "

$ws.Range("D6").Value = "0.6226"
$ws.Range("E6").Value = "  -5.23%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D10").Value = "0.2908"
$ws.Range("E10").Value = "  -0.56%  "

$ws.Range("D11").Value = "22.73"
$ws.Range("E11").Value = "  -1.03%  "

$ws.Range("D13").Value = "1.832.04"
$ws.Range("E13").Value = "  +0.43%  "

$ws.Range("D14").Value = "4.954"
$ws.Range("E14").Value = "  -0.50%  "

$ws.Range("D15").Value = "0.6639"
$ws.Range("E15").Value = "  -0.31%  "

$ws.Range("D16").Value = "82.28"
$ws.Range("E16").Value = "  -0.54%  "

$ws.Range("D17").Value = "0.000009122"
$ws.Range("E17").Value = "  +7.35%  "

$ws.Range("D18").Value = "5.985"
$ws.Range("E18").Value = "  -2.20%  "

$ws.Range("D19").Value = "28.919.11"
$ws.Range("E19").Value = "  -0.76%  "

$ws.Range("D20").Value = "224.56"
$ws.Range("E20").Value = "  -1.21%  "

$ws.Range("D21").Value = "12.32"
$ws.Range("E21").Value = "  -1.11%  "

$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.07%  "

$ws.Range("D23").Value = "7.183"
$ws.Range("E23").Value = "  +0.61%  "

$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  +0.11%  "

$ws.Range("D25").Value = "159.64"
$ws.Range("E25").Value = "  +0.68%  "

$ws.Range("D26").Value = "8.391"
$ws.Range("E26").Value = "  -2.35%  "

$ws.Range("D27").Value = "0.1354"
$ws.Range("E27").Value = "  -2.78%  "

$ws.Range("D28").Value = "17.80"
$ws.Range("E28").Value = "  -0.70%  "

$ws.Range("D29").Value = "1.495"
$ws.Range("E29").Value = "  -1.55%  "

$ws.Range("D30").Value = "4.028"
$ws.Range("E30").Value = "  -0.45%  "

$ws.Range("D31").Value = "4.045"
$ws.Range("E31").Value = "  -1.67%  "

$ws.Range("D32").Value = "1.201"
$ws.Range("E32").Value = "  +0.67%  "

$ws.Range("D33").Value = "0.05197"
$ws.Range("E33").Value = "  -1.13%  "

$ws.Range("D34").Value = "1.835"
$ws.Range("E34").Value = "  -1.64%  "

$ws.Range("D35").Value = "1.151"
$ws.Range("E35").Value = "  +0.59%  "

$ws.Range("D36").Value = "0.7300"
$ws.Range("E36").Value = "  -0.97%  "

$ws.Range("D37").Value = "2.611"
$ws.Range("E37").Value = "  -1.52%  "

$ws.Range("D38").Value = "1.278.24"
$ws.Range("E38").Value = "  -2.07%  "

$ws.Range("D39").Value = "2.763"
$ws.Range("E39").Value = "  +1.16%  "

$ws.Range("D40").Value = "0.01784"
$ws.Range("E40").Value = "  -0.72%  "

$ws.Range("D41").Value = "6.400"
$ws.Range("E41").Value = "  +7.27%  "

$ws.Range("D42").Value = "0.8898"
$ws.Range("E42").Value = "  -3.38%  "

$ws.Range("D43").Value = "1.000"
$ws.Range("E43").Value = "  +0.12%  "

$ws.Range("D44").Value = "101.37"
$ws.Range("E44").Value = "  -1.00%  "

$ws.Range("D45").Value = "1.979.76"
$ws.Range("E45").Value = "  +2.89%  "

$ws.Range("D46").Value = "0.5107"
$ws.Range("E46").Value = "  -0.65%  "

$ws.Range("D47").Value = "63.61"
$ws.Range("E47").Value = "  +0.63%  "

$ws.Range("E48").Value = "  -0.58%  "

$ws.Range("D49").Value = "0.07427"
$ws.Range("E49").Value = "  -15.39%  "

$ws.Range("D50").Value = "0.3969"
$ws.Range("E50").Value = "  -0.92%  "

$ws.Range("D51").Value = "8.910"
$ws.Range("E51").Value = "  +1.99%  "

# Restore the default style on the range (the temporary "@" text NumberFormat
# bumped the cell style index); this keeps the values as text while leaving
# the cell formatting/style identical to the original file.
$ws.Range("D2:E51").Style = "Normal"
